$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @("28.143.36", "  +3.85%  ")
    3  = @("1.785.47", "  +0.19%  ")
    4  = @("0.9939", "  -1.12%  ")
    5  = @("333.49", "  -1.09%  ")
    6  = @("0.9940", "  -0.81%  ")
    7  = @("0.3825", "  +0.31%  ")
    8  = @("0.3429", "  +0.60%  ")
    9  = @("47.92", "  -0.41%  ")
    10 = @("1.152", "  -2.89%  ")
    11 = @("0.07463", "  +0.37%  ")
    12 = @("22.95", "  +6.16%  ")
    13 = @("0.9907", "  -1.19%  ")
    14 = @("6.409", "  -0.60%  ")
    15 = @("1.780.20", "  +0.12%  ")
    16 = @("7.142", "  +1.09%  ")
    17 = @("0.00001083", "  -0.08%  ")
    18 = @("0.06610", "  -0.37%  ")
    19 = @("83.02", "  -0.36%  ")
    20 = @("0.9952", "  -0.69%  ")
    21 = @("17.56", "  +1.43%  ")
    22 = @("6.452", "  -1.23%  ")
    23 = @("28.102.03", "  +3.69%  ")
    24 = @("12.15", "  -0.45%  ")
    25 = @("2.375", "  +0.20%  ")
    26 = @("1.452", "  -0.03%  ")
    27 = @("20.94", "  -0.49%  ")
    28 = @("2.464", "  -1.51%  ")
    29 = @("154.28", $null)
    30 = @("1.979.38", "  +0.01%  ")
    31 = @("134.70", "  +0.61%  ")
    32 = @("6.210", "  +3.13%  ")
    33 = @("3.946", "  -0.91%  ")
    34 = @("0.08798", "  +1.44%  ")
    35 = @("12.87", "  -1.61%  ")
    36 = @("0.02443", "  +5.56%  ")
    37 = @("0.6922", "  +1.42%  ")
    38 = @("5.363", "  -0.38%  ")
    39 = @("0.06385", "  +1.72%  ")
    40 = @("0.2200", "  +1.30%  ")
    41 = @("1.512", "  -7.08%  ")
    42 = @("1.241", "  +0.35%  ")
    43 = @("8.403", "  -1.55%  ")
    44 = @("14.24", "  +0.33%  ")
    45 = @("0.9938", "  -0.75%  ")
    46 = @("0.6351", "  -1.20%  ")
    47 = @("3.825", "  -0.77%  ")
    48 = @("132.80", "  +1.10%  ")
    49 = @("2.100", "  -0.93%  ")
    50 = @("0.07417", "  +4.72%  ")
    51 = @("78.82", "  +0.16%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $dVal

    if ($eVal -ne $null) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $eVal
    }
}
